$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing B-column values for revised historical/forecast netback figures ---
$updates = @{
    171 = 561.92669584312637
    173 = 848.92898689293725
    176 = 728.23648918640038
    177 = 707.53396319493686
    178 = 641.48306145840672
    179 = 563.3071678513752
    180 = 537.46710932784106
    181 = 464.20843872713505
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

# --- Append new monthly rows (2023 manual forecast) ---
$newRows = @(
    @(182, 44927, 397.23604124929227, 2),
    @(183, 44958, 397.34038136656, 1),
    @(184, 44986, 397.84397137970029, 2),
    @(185, 45017, 398.16037878233846, 1),
    @(186, 45047, 399.68264381252743, 2),
    @(187, 45078, 399.25312666369229, 1),
    @(188, 45108, 397.33526810467691, 2),
    @(189, 45139, 397.00785842393844, 1),
    @(190, 45170, 396.70599348008614, 2),
    @(191, 45200, 397.67744701382151, 1),
    @(192, 45231, 397.12623831080612, 2),
    @(193, 45261, 396.19664348308919, 1),
)

foreach ($r in $newRows) {
    $row = $r[0]
    $dateSerial = $r[1]
    $value = $r[2]
    $styleFlag = $r[3]

    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.Value = $dateSerial
    $dateCell.NumberFormat = "mmm-yy"
    if ($styleFlag -eq 2) {
        $dateCell.Interior.Color = 65535
    }

    $ws.Cells.Item($row, 2).Value = $value
}

# --- Update view selection to match the active editing position ---
$ws.Range("B158:B193").Select()
